# Actualización automática 2025-08-08 16:15:10
#
# LOZANO MOLINA TITO / RENOVA&DISEÑA S.A. registró una venta adicional de
# 167.45 en el grupo PORCELANATO durante el mes de agosto. Se propaga el
# valor a las tres hojas del libro (VENTAS POR GRUPO, VENTA MENSUAL y
# CUMPLIMIENTO MENSUAL) y se ensancha ligeramente la columna de agosto.

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO": venta del cliente por grupo de producto ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Cliente RENOVA&DISEÑA S.A. (fila 19), columna M = PORCELANATO
$wsGrupo.Range("M19").Value = 167.45

# Fila resumen 29: cantidad de asesores/clientes con venta en el grupo
$wsGrupo.Range("M29").Value = "1 de 27"

# --- Hoja "VENTA MENSUAL": venta del cliente por mes ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Columna F = agosto, fila 19 = RENOVA&DISEÑA S.A.
$wsMensual.Range("F19").Value = 1837.88

# Fila 29 = total de agosto
$wsMensual.Range("F29").Value = 1837.88

# --- Hoja "CUMPLIMIENTO MENSUAL": cumplimiento de presupuesto por grupo ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Ensanchar la columna CUMPLIMIENTO (columna F, la 6ta)
# (25.17 aterriza exactamente en el ancho almacenado de 26 tras la
# cuantización a la cuadrícula de píxeles que aplica Excel/COM)
$wsCumplimiento.Columns.Item(6).ColumnWidth = 25.17

# Fila 16 = grupo PORCELANATO
$wsCumplimiento.Range("D16").Value = 167.45
$wsCumplimiento.Range("E16").Value = 23737.13
$wsCumplimiento.Range("F16").Value = 0.007004933782563842

# Fila 19 = TOTAL
$wsCumplimiento.Range("D19").Value = 1837.88
$wsCumplimiento.Range("E19").Value = 35662.12093005038
$wsCumplimiento.Range("F19").Value = 0.04901013211781621
